# Add new row 16 to Sheet1: "Sliding Window" / "3. Longest Substring Without
# Repeating Characters" / solution summary (with a bold inline phrase).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Copy the formatting of an existing "category" row (row 5 has the same
#    A=border/wrap, B=yellow-fill, C=border/wrap style pattern, and even the
#    same row height of 43.2, that the new row ends up with) down onto row 16.
$srcRow = $ws.Range("A5:C5")
$dstRow = $ws.Range("A16:C16")
$srcRow.Copy()
$dstRow.PasteSpecial(-4122)  # xlPasteFormats
$ws.Rows.Item(16).RowHeight = 43.2

# 2. Fill in the plain-text cells.
$ws.Range("A16").Value = "Sliding Window"
$ws.Range("B16").Value = "3. Longest Substring Without Repeating Characters"

# 3. Fill in the rich-text solution summary in C16, then bold the phrase
#    "before this" in the middle of it.
$prefix = "The idea is to use a hashset & keep track of the start of the window using `"winStart`", iter over the chars in array, if currChar isnt present in hs then add it, "
$bolded = "before this"
$suffix = " put a while(hs.contains(currChar)) loop do hs.remove( charArr[winstart] ) then winStart++, since we are sliding the window to the right if we detect a char same as currChar in hashset"
$full = $prefix + $bolded + $suffix

$ws.Range("C16").Value = $full

$boldStart = $prefix.Length + 1
$boldLen = $bolded.Length
$ws.Range("C16").Characters($boldStart, $boldLen).Font.Bold = $true

$suffixStart = $boldStart + $boldLen
$suffixLen = $full.Length - $suffixStart + 1
$ws.Range("C16").Characters($suffixStart, $suffixLen).Font.Size = 11

# 4. Match the saved selection/active cell.
$ws.Range("B16").Select() | Out-Null
